$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 3
    4  = 3.5
    5  = 3.5
    6  = 5.5
    7  = 3.5
    8  = 7
    9  = 5
    10 = 19.5
    11 = 1.5
    12 = 19.5
    13 = 2.5
    14 = 5.5
    16 = 6.5
    18 = 12.5
    19 = 4
    20 = 22.5
    21 = 2.5
    22 = 15
    23 = 5.5
    24 = 41.5
    25 = 7.5
    26 = 17.5
    27 = 2
    28 = 18.5
    29 = 2.5
    30 = 12
    31 = 3.5
    32 = 8.5
    33 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
